# Updates odds/market data cells on Sheet1 to match the 2026-01-27 refresh.
# Values below were taken from the authoritative diff of the commit
# "Atualizando o arquivo XLSX" (231 cell value updates across rows 2-25).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 8).Value = 1.94  # H2: 1.89 -> 1.94
$ws.Cells.Item(2, 9).Value = 160  # I2: 1000 -> 160
$ws.Cells.Item(2, 10).Value = 3  # J2: 1.37 -> 3
$ws.Cells.Item(2, 11).Value = 4.5  # K2: 4.8 -> 4.5
$ws.Cells.Item(2, 12).Value = 1.3  # L2: 1.01 -> 1.3
$ws.Cells.Item(2, 13).Value = 1.06  # M2: 1.05 -> 1.06
$ws.Cells.Item(2, 14).Value = 1.7  # N2: 1.08 -> 1.7
$ws.Cells.Item(2, 15).Value = 1.06  # O2: 1.05 -> 1.06
$ws.Cells.Item(2, 16).Value = 1.7  # P2: 1.08 -> 1.7
$ws.Cells.Item(2, 17).Value = 1.65  # Q2: 1.05 -> 1.65
$ws.Cells.Item(2, 19).Value = 1.65  # S2: 1.05 -> 1.65
$ws.Cells.Item(2, 22).Value = 1.02  # V2: 1.01 -> 1.02
# Row 3
$ws.Cells.Item(3, 6).Value = 3.4  # F3: 3.5 -> 3.4
$ws.Cells.Item(3, 7).Value = 4.7  # G3: 4.9 -> 4.7
$ws.Cells.Item(3, 9).Value = 2.1  # I3: 2.3 -> 2.1
$ws.Cells.Item(3, 10).Value = 3.4  # J3: 3.05 -> 3.4
$ws.Cells.Item(3, 11).Value = 4.5  # K3: 3.95 -> 4.5
$ws.Cells.Item(3, 12).Value = 1.31  # L3: 1.41 -> 1.31
$ws.Cells.Item(3, 13).Value = 1.05  # M3: 1.07 -> 1.05
$ws.Cells.Item(3, 14).Value = 4.4  # N3: 3.2 -> 4.4
$ws.Cells.Item(3, 15).Value = 1.22  # O3: 1.33 -> 1.22
$ws.Cells.Item(3, 16).Value = 2.2  # P3: 1.84 -> 2.2
$ws.Cells.Item(3, 17).Value = 1.66  # Q3: 1.94 -> 1.66
$ws.Cells.Item(3, 18).Value = 1.49  # R3: 1.32 -> 1.49
$ws.Cells.Item(3, 19).Value = 2.48  # S3: 3.15 -> 2.48
$ws.Cells.Item(3, 20).Value = 1.62  # T3: 1.79 -> 1.62
$ws.Cells.Item(3, 21).Value = 2.3  # U3: 2.02 -> 2.3
$ws.Cells.Item(3, 22).Value = 1.9  # V3: 1.83 -> 1.9
$ws.Cells.Item(3, 23).Value = 1.29  # W3: 1.26 -> 1.29
$ws.Cells.Item(3, 24).Value = 24  # X3: 16 -> 24
$ws.Cells.Item(3, 25).Value = 14  # Y3: 10.5 -> 14
$ws.Cells.Item(3, 26).Value = 16  # Z3: 15 -> 16
$ws.Cells.Item(3, 27).Value = 27  # AA3: 29 -> 27
$ws.Cells.Item(3, 28).Value = 19.5  # AB3: 17 -> 19.5
$ws.Cells.Item(3, 29).Value = 10.5  # AC3: 9.4 -> 10.5
$ws.Cells.Item(3, 31).Value = 22  # AE3: 26 -> 22
$ws.Cells.Item(3, 33).Value = 19  # AG3: 19.5 -> 19
$ws.Cells.Item(3, 34).Value = 18.5  # AH3: 22 -> 18.5
$ws.Cells.Item(3, 35).Value = 34  # AI3: 44 -> 34
$ws.Cells.Item(3, 36).Value = 85  # AJ3: 100 -> 85
$ws.Cells.Item(3, 37).Value = 48  # AK3: 60 -> 48
$ws.Cells.Item(3, 38).Value = 50  # AL3: 70 -> 50
$ws.Cells.Item(3, 39).Value = 75  # AM3: 120 -> 75
$ws.Cells.Item(3, 40).Value = 36  # AN3: 70 -> 36
$ws.Cells.Item(3, 41).Value = 11.5  # AO3: 19 -> 11.5
# Row 4
$ws.Cells.Item(4, 6).Value = 3.3  # F4: 3.2 -> 3.3
$ws.Cells.Item(4, 8).Value = 2.16  # H4: 2.14 -> 2.16
$ws.Cells.Item(4, 9).Value = 2.52  # I4: 2.66 -> 2.52
$ws.Cells.Item(4, 19).Value = 3.75  # S4: 3.8 -> 3.75
$ws.Cells.Item(4, 20).Value = 1.83  # T4: 1.81 -> 1.83
$ws.Cells.Item(4, 23).Value = 1.32  # W4: 1.31 -> 1.32
# Row 5
$ws.Cells.Item(5, 6).Value = 1.95  # F5: 1.94 -> 1.95
$ws.Cells.Item(5, 7).Value = 2.16  # G5: 2.14 -> 2.16
$ws.Cells.Item(5, 8).Value = 3.85  # H5: 3.9 -> 3.85
$ws.Cells.Item(5, 12).Value = 1.42  # L5: 1.41 -> 1.42
$ws.Cells.Item(5, 23).Value = 1.86  # W5: 1.87 -> 1.86
$ws.Cells.Item(5, 25).Value = 15.5  # Y5: 17.5 -> 15.5
$ws.Cells.Item(5, 26).Value = 970  # Z5: 34 -> 970
$ws.Cells.Item(5, 30).Value = 18  # AD5: 21 -> 18
$ws.Cells.Item(5, 32).Value = 14.5  # AF5: 15 -> 14.5
$ws.Cells.Item(5, 33).Value = 11.5  # AG5: 13 -> 11.5
$ws.Cells.Item(5, 35).Value = 70  # AI5: 80 -> 70
# Row 6
$ws.Cells.Item(6, 8).Value = 21  # H6: 10 -> 21
# Row 8
$ws.Cells.Item(8, 8).Value = 1.09  # H8: 1.04 -> 1.09
$ws.Cells.Item(8, 15).Value = 1.66  # O8: 1.63 -> 1.66
$ws.Cells.Item(8, 19).Value = 1.67  # S8: 1.63 -> 1.67
$ws.Cells.Item(8, 22).Value = 1.17  # V8: 1.05 -> 1.17
# Row 9
$ws.Cells.Item(9, 11).Value = 4  # K9: 4.1 -> 4
$ws.Cells.Item(9, 14).Value = 2.84  # N9: 2.86 -> 2.84
$ws.Cells.Item(9, 16).Value = 1.61  # P9: 1.63 -> 1.61
$ws.Cells.Item(9, 23).Value = 2.52  # W9: 2.56 -> 2.52
# Row 10
$ws.Cells.Item(10, 7).Value = 3.7  # G10: 3.75 -> 3.7
$ws.Cells.Item(10, 9).Value = 2.16  # I10: 2.18 -> 2.16
$ws.Cells.Item(10, 10).Value = 3.85  # J10: 3.8 -> 3.85
$ws.Cells.Item(10, 16).Value = 2.52  # P10: 2.5 -> 2.52
$ws.Cells.Item(10, 22).Value = 1.86  # V10: 1.85 -> 1.86
$ws.Cells.Item(10, 33).Value = 19  # AG10: 18.5 -> 19
# Row 11
$ws.Cells.Item(11, 6).Value = 4.5  # F11: 4.6 -> 4.5
$ws.Cells.Item(11, 8).Value = 1.89  # H11: 1.87 -> 1.89
$ws.Cells.Item(11, 9).Value = 1.9  # I11: 1.89 -> 1.9
$ws.Cells.Item(11, 16).Value = 2.14  # P11: 2.16 -> 2.14
$ws.Cells.Item(11, 17).Value = 1.84  # Q11: 1.82 -> 1.84
$ws.Cells.Item(11, 19).Value = 3.15  # S11: 3.1 -> 3.15
$ws.Cells.Item(11, 20).Value = 1.77  # T11: 1.78 -> 1.77
$ws.Cells.Item(11, 22).Value = 2.1  # V11: 2.12 -> 2.1
$ws.Cells.Item(11, 25).Value = 10  # Y11: 10.5 -> 10
$ws.Cells.Item(11, 27).Value = 21  # AA11: 20 -> 21
$ws.Cells.Item(11, 28).Value = 18  # AB11: 18.5 -> 18
# Row 12
$ws.Cells.Item(12, 6).Value = 3.6  # F12: 3.55 -> 3.6
$ws.Cells.Item(12, 7).Value = 3.65  # G12: 3.6 -> 3.65
$ws.Cells.Item(12, 8).Value = 2.12  # H12: 2.16 -> 2.12
$ws.Cells.Item(12, 9).Value = 2.16  # I12: 2.2 -> 2.16
$ws.Cells.Item(12, 10).Value = 3.85  # J12: 3.8 -> 3.85
$ws.Cells.Item(12, 12).Value = 1.32  # L12: 1.31 -> 1.32
$ws.Cells.Item(12, 13).Value = 1.05  # M12: 1.04 -> 1.05
$ws.Cells.Item(12, 14).Value = 5.4  # N12: 5.5 -> 5.4
$ws.Cells.Item(12, 15).Value = 1.22  # O12: 1.21 -> 1.22
$ws.Cells.Item(12, 16).Value = 2.42  # P12: 2.46 -> 2.42
$ws.Cells.Item(12, 17).Value = 1.68  # Q12: 1.66 -> 1.68
$ws.Cells.Item(12, 18).Value = 1.58  # R12: 1.59 -> 1.58
$ws.Cells.Item(12, 19).Value = 2.66  # S12: 2.62 -> 2.66
$ws.Cells.Item(12, 20).Value = 1.61  # T12: 1.59 -> 1.61
$ws.Cells.Item(12, 21).Value = 2.56  # U12: 2.66 -> 2.56
$ws.Cells.Item(12, 23).Value = 1.37  # W12: 1.38 -> 1.37
$ws.Cells.Item(12, 25).Value = 13  # Y12: 13.5 -> 13
$ws.Cells.Item(12, 26).Value = 15.5  # Z12: 16 -> 15.5
$ws.Cells.Item(12, 28).Value = 18  # AB12: 18.5 -> 18
$ws.Cells.Item(12, 30).Value = 11  # AD12: 10.5 -> 11
$ws.Cells.Item(12, 38).Value = 40  # AL12: 38 -> 40
$ws.Cells.Item(12, 39).Value = 65  # AM12: 60 -> 65
$ws.Cells.Item(12, 40).Value = 26  # AN12: 25 -> 26
$ws.Cells.Item(12, 41).Value = 11.5  # AO12: 11 -> 11.5
# Row 14
$ws.Cells.Item(14, 10).Value = 3.45  # J14: 3.6 -> 3.45
$ws.Cells.Item(14, 11).Value = 3.85  # K14: 3.95 -> 3.85
# Row 15
$ws.Cells.Item(15, 6).Value = 2.36  # F15: 2.54 -> 2.36
$ws.Cells.Item(15, 7).Value = 2.72  # G15: 2.7 -> 2.72
$ws.Cells.Item(15, 8).Value = 2.74  # H15: 2.72 -> 2.74
$ws.Cells.Item(15, 9).Value = 3.2  # I15: 2.8 -> 3.2
$ws.Cells.Item(15, 10).Value = 3.45  # J15: 3.65 -> 3.45
$ws.Cells.Item(15, 22).Value = 1.45  # V15: 1.55 -> 1.45
$ws.Cells.Item(15, 23).Value = 1.58  # W15: 1.59 -> 1.58
$ws.Cells.Item(15, 24).Value = 19  # X15: 19.5 -> 19
$ws.Cells.Item(15, 25).Value = 14.5  # Y15: 15 -> 14.5
$ws.Cells.Item(15, 28).Value = 13.5  # AB15: 14 -> 13.5
$ws.Cells.Item(15, 32).Value = 19  # AF15: 19.5 -> 19
$ws.Cells.Item(15, 36).Value = 38  # AJ15: 40 -> 38
# Row 16
$ws.Cells.Item(16, 6).Value = 2.32  # F16: 2.3 -> 2.32
$ws.Cells.Item(16, 10).Value = 2.92  # J16: 2.9 -> 2.92
$ws.Cells.Item(16, 12).Value = 1.65  # L16: 1.66 -> 1.65
$ws.Cells.Item(16, 14).Value = 2.38  # N16: 2.34 -> 2.38
$ws.Cells.Item(16, 15).Value = 1.65  # O16: 1.66 -> 1.65
$ws.Cells.Item(16, 16).Value = 1.46  # P16: 1.44 -> 1.46
$ws.Cells.Item(16, 18).Value = 1.16  # R16: 1.15 -> 1.16
$ws.Cells.Item(16, 21).Value = 1.67  # U16: 1.66 -> 1.67
$ws.Cells.Item(16, 24).Value = 7.8  # X16: 8.800000000000001 -> 7.8
$ws.Cells.Item(16, 28).Value = 6.8  # AB16: 7.2 -> 6.8
$ws.Cells.Item(16, 29).Value = 8.4  # AC16: 8.6 -> 8.4
$ws.Cells.Item(16, 31).Value = 75  # AE16: 80 -> 75
$ws.Cells.Item(16, 33).Value = 13  # AG16: 15.5 -> 13
$ws.Cells.Item(16, 37).Value = 42  # AK16: 46 -> 42
# Row 17
$ws.Cells.Item(17, 7).Value = 3.8  # G17: 3.85 -> 3.8
$ws.Cells.Item(17, 8).Value = 2.48  # H17: 2.46 -> 2.48
$ws.Cells.Item(17, 12).Value = 1.51  # L17: 1.57 -> 1.51
$ws.Cells.Item(17, 14).Value = 2.52  # N17: 2.5 -> 2.52
$ws.Cells.Item(17, 25).Value = 970  # Y17: 8.199999999999999 -> 970
$ws.Cells.Item(17, 29).Value = 7.4  # AC17: 970 -> 7.4
# Row 19
$ws.Cells.Item(19, 11).Value = 2.98  # K19: 3 -> 2.98
$ws.Cells.Item(19, 17).Value = 3.4  # Q19: 3.35 -> 3.4
$ws.Cells.Item(19, 22).Value = 1.22  # V19: 1.21 -> 1.22
$ws.Cells.Item(19, 28).Value = 5.8  # AB19: 6 -> 5.8
$ws.Cells.Item(19, 31).Value = 160  # AE19: 130 -> 160
$ws.Cells.Item(19, 33).Value = 13.5  # AG19: 15 -> 13.5
$ws.Cells.Item(19, 38).Value = 120  # AL19: 95 -> 120
# Row 20
$ws.Cells.Item(20, 6).Value = 2.84  # F20: 2.9 -> 2.84
$ws.Cells.Item(20, 7).Value = 3.05  # G20: 3.1 -> 3.05
$ws.Cells.Item(20, 8).Value = 3.15  # H20: 3.1 -> 3.15
$ws.Cells.Item(20, 9).Value = 3.5  # I20: 3.4 -> 3.5
$ws.Cells.Item(20, 10).Value = 2.7  # J20: 2.68 -> 2.7
$ws.Cells.Item(20, 11).Value = 2.9  # K20: 2.84 -> 2.9
$ws.Cells.Item(20, 14).Value = 2.02  # N20: 2.08 -> 2.02
$ws.Cells.Item(20, 15).Value = 1.83  # O20: 1.81 -> 1.83
$ws.Cells.Item(20, 16).Value = 1.34  # P20: 1.35 -> 1.34
$ws.Cells.Item(20, 17).Value = 3.55  # Q20: 3.5 -> 3.55
$ws.Cells.Item(20, 20).Value = 2.6  # T20: 2.5 -> 2.6
$ws.Cells.Item(20, 21).Value = 1.48  # U20: 1.57 -> 1.48
$ws.Cells.Item(20, 23).Value = 1.48  # W20: 1.47 -> 1.48
$ws.Cells.Item(20, 31).Value = 1000  # AE20: 70 -> 1000
$ws.Cells.Item(20, 33).Value = 18.5  # AG20: 970 -> 18.5
# Row 22
$ws.Cells.Item(22, 7).Value = 2.24  # G22: 2.18 -> 2.24
$ws.Cells.Item(22, 8).Value = 4.5  # H22: 4.6 -> 4.5
$ws.Cells.Item(22, 9).Value = 5.3  # I22: 5.2 -> 5.3
$ws.Cells.Item(22, 10).Value = 2.94  # J22: 3 -> 2.94
$ws.Cells.Item(22, 12).Value = 1.54  # L22: 1.51 -> 1.54
$ws.Cells.Item(22, 13).Value = 1.12  # M22: 1.09 -> 1.12
$ws.Cells.Item(22, 14).Value = 2.58  # N22: 2.72 -> 2.58
$ws.Cells.Item(22, 15).Value = 1.51  # O22: 1.47 -> 1.51
$ws.Cells.Item(22, 16).Value = 1.51  # P22: 1.57 -> 1.51
$ws.Cells.Item(22, 17).Value = 2.6  # Q22: 2.46 -> 2.6
$ws.Cells.Item(22, 18).Value = 1.19  # R22: 1.21 -> 1.19
$ws.Cells.Item(22, 19).Value = 5  # S22: 4.8 -> 5
$ws.Cells.Item(22, 20).Value = 2.12  # T22: 2.06 -> 2.12
$ws.Cells.Item(22, 21).Value = 1.74  # U22: 1.78 -> 1.74
$ws.Cells.Item(22, 22).Value = 1.25  # V22: 1.26 -> 1.25
$ws.Cells.Item(22, 23).Value = 1.82  # W22: 1.84 -> 1.82
$ws.Cells.Item(22, 24).Value = 9.199999999999999  # X22: 11.5 -> 9.199999999999999
$ws.Cells.Item(22, 25).Value = 13  # Y22: 14.5 -> 13
$ws.Cells.Item(22, 29).Value = 7.4  # AC22: 7.6 -> 7.4
$ws.Cells.Item(22, 33).Value = 12.5  # AG22: 11.5 -> 12.5
$ws.Cells.Item(22, 38).Value = 60  # AL22: 55 -> 60
$ws.Cells.Item(22, 39).Value = 230  # AM22: 190 -> 230
$ws.Cells.Item(22, 40).Value = 28  # AN22: 26 -> 28
$ws.Cells.Item(22, 41).Value = 170  # AO22: 130 -> 170
# Row 23
$ws.Cells.Item(23, 6).Value = 2.34  # F23: 2.38 -> 2.34
$ws.Cells.Item(23, 7).Value = 2.5  # G23: 2.58 -> 2.5
$ws.Cells.Item(23, 8).Value = 3.85  # H23: 3.65 -> 3.85
$ws.Cells.Item(23, 9).Value = 4.3  # I23: 4.2 -> 4.3
$ws.Cells.Item(23, 11).Value = 2.98  # K23: 3 -> 2.98
$ws.Cells.Item(23, 14).Value = 2.36  # N23: 2.4 -> 2.36
$ws.Cells.Item(23, 15).Value = 1.62  # O23: 1.61 -> 1.62
$ws.Cells.Item(23, 16).Value = 1.44  # P23: 1.46 -> 1.44
$ws.Cells.Item(23, 23).Value = 1.66  # W23: 1.64 -> 1.66
# Row 24
$ws.Cells.Item(24, 6).Value = 3.45  # F24: 3.55 -> 3.45
$ws.Cells.Item(24, 7).Value = 4.7  # G24: 3.95 -> 4.7
$ws.Cells.Item(24, 8).Value = 2.4  # H24: 2.46 -> 2.4
$ws.Cells.Item(24, 9).Value = 2.72  # I24: 2.64 -> 2.72
$ws.Cells.Item(24, 10).Value = 2.58  # J24: 2.86 -> 2.58
$ws.Cells.Item(24, 12).Value = 1.01  # L24: 1.62 -> 1.01
$ws.Cells.Item(24, 13).Value = 1.01  # M24: 1.12 -> 1.01
$ws.Cells.Item(24, 14).Value = 1.25  # N24: 2.44 -> 1.25
$ws.Cells.Item(24, 16).Value = 1.25  # P24: 1.48 -> 1.25
$ws.Cells.Item(24, 18).Value = 1.12  # R24: 1.17 -> 1.12
$ws.Cells.Item(24, 19).Value = 5.3  # S24: 6 -> 5.3
$ws.Cells.Item(24, 20).Value = 2.28  # T24: 2.18 -> 2.28
$ws.Cells.Item(24, 21).Value = 1.68  # U24: 1.72 -> 1.68
$ws.Cells.Item(24, 22).Value = 1.6  # V24: 1.61 -> 1.6
$ws.Cells.Item(24, 23).Value = 1.33  # W24: 1.34 -> 1.33
$ws.Cells.Item(24, 24).Value = 9  # X24: 9.4 -> 9
# Row 25
$ws.Cells.Item(25, 6).Value = 2.1  # F25: 2.08 -> 2.1
$ws.Cells.Item(25, 10).Value = 3  # J25: 3.05 -> 3
$ws.Cells.Item(25, 13).Value = 1.12  # M25: 1.11 -> 1.12
$ws.Cells.Item(25, 14).Value = 2.7  # N25: 2.8 -> 2.7
$ws.Cells.Item(25, 15).Value = 1.48  # O25: 1.49 -> 1.48
$ws.Cells.Item(25, 16).Value = 1.56  # P25: 1.59 -> 1.56
$ws.Cells.Item(25, 17).Value = 2.54  # Q25: 2.46 -> 2.54
$ws.Cells.Item(25, 18).Value = 1.2  # R25: 1.21 -> 1.2
$ws.Cells.Item(25, 19).Value = 5.1  # S25: 4.9 -> 5.1
$ws.Cells.Item(25, 20).Value = 2.08  # T25: 2.06 -> 2.08
$ws.Cells.Item(25, 21).Value = 1.78  # U25: 1.79 -> 1.78
$ws.Cells.Item(25, 28).Value = 7.2  # AB25: 7.4 -> 7.2
$ws.Cells.Item(25, 37).Value = 38  # AK25: 29 -> 38
$ws.Cells.Item(25, 38).Value = 70  # AL25: 55 -> 70
